$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New values (this handback run regenerated the two e2e/*.md files under new
# GUIDs and re-ran the xliff round trip, producing new timestamps).
# ---------------------------------------------------------------------------
$oldGuid1 = "4930cd6e-3d3b-4858-a9c5-80fc33b1aac0"
$newGuid1 = "9a023f4c-dea3-4eab-91ea-79f72a209048"
$oldGuid2 = "c4afc7f6-4f30-4781-b948-86c6556fe580"
$newGuid2 = "ffffefd082a9-af3c-47af-88e1-af967a44b92d"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"
$newPath1 = "e2e\$newGuid1.md"
$newPath2 = "e2e\$newGuid2.md"

$newXlf1ZhCn = "$newGuid1.074d7c1959bb106be22360d7b6cb090df0a970a1.zh-cn.xlf"
$newXlf1DeDe = "$newGuid1.074d7c1959bb106be22360d7b6cb090df0a970a1.de-de.xlf"

$newHandoffDate = "2016-08-20 13:03:55"
$newZhCnHandoffDate = "2016-08-20 13:03:51"
$newZhCnHandbackDate = "2016-08-20 13:04:12"
$newDeDeHandbackDate = "2016-08-20 13:04:19"

# Original external targets (unchanged relationships) keyed by sheet + GUID.
$overviewUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid1.md"
$overviewUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid2.md"

$zhcnUrlA1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid1.md"
$zhcnUrlI1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/46f76c0f64146d60e44e17eb05190cf24ad09679/e2e/$oldGuid1.md"
$zhcnUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid2.md"
$zhcnUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/46f76c0f64146d60e44e17eb05190cf24ad09679/e2e/$oldGuid2.md"

$dedeUrlA1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid1.md"
$dedeUrlI1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1a037df79c77af73750130ab13621fe7c4a14f04/e2e/$oldGuid1.md"
$dedeUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/33aba1a97c9ef7eda9d503ed70dc27a3cb2f4d74/e2e/$oldGuid2.md"
$dedeUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1a037df79c77af73750130ab13621fe7c4a14f04/e2e/$oldGuid2.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("G2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("G3").Value = $newHandoffDate

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewUrl1, "", "", $newPath1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $overviewUrl2, "", "", $newPath2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $newXlf1ZhCn
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $newXlf1ZhCn
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $newXlf1ZhCn
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $newXlf1ZhCn
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnUrlA1, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhcnUrlI1, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhcnUrlA2, "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhcnUrlI2, "", "", $newFile2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $newXlf1DeDe
$wsDeDe.Range("H2").Value = $newHandoffDate
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $newXlf1DeDe
$wsDeDe.Range("K2").Value = $newDeDeHandbackDate

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $newXlf1DeDe
$wsDeDe.Range("H3").Value = $newHandoffDate
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $newXlf1DeDe
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeUrlA1, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $dedeUrlI1, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $dedeUrlA2, "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $dedeUrlI2, "", "", $newFile2)

Write-Host "Handback status report regenerated."
